# Update gh-pages output data: 想去人数 (interest count) bumps on two sheets
# that share identical underlying data ("展览" sheet 1 and "全部类型" sheet 4).
#   F2: 303 -> 305
#   F4: 172 -> 173

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 305
    $ws.Range("F4").Value = 173
}
